$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.790495038032532
$ws.Range("B1").Value = 4.144339084625244
$ws.Range("C1").Value = 1.477611064910889
$ws.Range("D1").Value = 0.858174741268158
$ws.Range("E1").Value = 0.4637786448001862
